$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This string is shared by the Overview sheet (B/C columns) and the
#    per-locale sheets' Status column (C). Updating every cell that shows
#    this text reproduces the effect of the shared string being edited.
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Populate the handback report columns (F = Latest Target File,
#    G = Latest Handback File) for rows 2 and 3 of each locale sheet, and
#    refresh the Latest Handback DateTime (column H).
# ---------------------------------------------------------------------------
$zhTargetFile = "945cee0c-9390-4a37-8313-f04094598067.md"
$zhHandbackFile = "945cee0c-9390-4a37-8313-f04094598067.b3728ad0985781609f24166094898695f8a4a3a0.zh-cn.xlf"
$zhTargetFileUrl = "https://github.com/OpenLocalizationTest/oltest/blob/3d824794cc69aa6b39265ac38dad1ca82ff67a28/e2e/945cee0c-9390-4a37-8313-f04094598067.md"
$zhHandbackFileUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5ae5b9ddc15d94936757ca021a29d1a2b7b3c7bc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/945cee0c-9390-4a37-8313-f04094598067.b3728ad0985781609f24166094898695f8a4a3a0.zh-cn.xlf"

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhTargetFileUrl, "", "", $zhTargetFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhHandbackFileUrl, "", "", $zhHandbackFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhTargetFileUrl, "", "", $zhTargetFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhHandbackFileUrl, "", "", $zhHandbackFile) | Out-Null

$wsZh.Range("F2:G3").Style = "HyperLink"
$wsZh.Range("H2").Value = "2016-03-12 22:55:45"
$wsZh.Range("H3").Value = "2016-03-12 22:55:45"

$deTargetFile = "945cee0c-9390-4a37-8313-f04094598067.md"
$deHandbackFile = "945cee0c-9390-4a37-8313-f04094598067.b3728ad0985781609f24166094898695f8a4a3a0.de-de.xlf"
$deTargetFileUrl = "https://github.com/OpenLocalizationTest/oltest/blob/3d824794cc69aa6b39265ac38dad1ca82ff67a28/e2e/945cee0c-9390-4a37-8313-f04094598067.md"
$deHandbackFileUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/55596834ec4af2bf1b7113501d77b7cf082e4b09/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/945cee0c-9390-4a37-8313-f04094598067.b3728ad0985781609f24166094898695f8a4a3a0.de-de.xlf"

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deTargetFileUrl, "", "", $deTargetFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deHandbackFileUrl, "", "", $deHandbackFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deTargetFileUrl, "", "", $deTargetFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deHandbackFileUrl, "", "", $deHandbackFile) | Out-Null

$wsDe.Range("F2:G3").Style = "HyperLink"
$wsDe.Range("H2").Value = "2016-03-12 22:55:51"
$wsDe.Range("H3").Value = "2016-03-12 22:55:51"
